$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.178967237472534
$ws.Range("B1").Value = 2.619049549102783
$ws.Range("C1").Value = 2.748536825180054
$ws.Range("D1").Value = 2.803977966308594
$ws.Range("E1").Value = 0.8100827932357788
